$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Cases Produced" row (row 3) with new solver values.
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 200
$ws.Range("E3").Value = 800
$ws.Range("F3").Value = 400
$ws.Range("G3").Value = 500
$ws.Range("H3").Value = 600

$excel.CalculateFullRebuild()
